$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (column C) from 45221 (2023-10-22)
# to 45224 (2023-10-25) for all data rows (rows 2-15).
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
